$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Plain/default style, used below to stamp cells back to the sheet's usual
# (unformatted) look after a write nudges them onto a new style index.
$normalStyle = $ws.Range("A2").Style

# Row 20 previously had several blank cells (B..K and M). The sheet's
# convention for "no data" is the literal text "nan", so backfill those.
"B20","C20","D20","E20","F20","G20","H20","I20","J20","K20","M20" | ForEach-Object {
    $ws.Range($_).Value = "nan"
}

# Add the new service event as row 21: card number 20, change date
# 2\2\2024, a correction note, and who serviced it. B21..K21 and M21 are
# left blank, same as the equivalent columns were before row 20 got its
# "nan" placeholders above.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "20"
$ws.Range("A21").Style = $normalStyle

"B21","C21","D21","E21","F21","G21","H21","I21","J21","K21","M21" | ForEach-Object {
    $ws.Range($_).Style = $normalStyle
}

$ws.Range("L21").Value = "2\2\2024"
$ws.Range("N21").Value = "تم تغيير 2 بليه النازع (يمين _ شمال)"
$ws.Range("O21").Value = "الخبير"
